$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.355.43"
$ws.Range("E2").Value = "  -0.72%  "
$ws.Range("D3").Value = "2.595.21"
$ws.Range("E3").Value = "  -1.71%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'510.18"
$ws.Range("E5").Value = "  -0.68%  "
$ws.Range("D6").Value = "'154.37"
$ws.Range("E6").Value = "  -2.87%  "
$ws.Range("E7").Value = "  +0.42%  "
$ws.Range("D8").Value = "'0.589"
$ws.Range("E8").Value = "  -4.32%  "
$ws.Range("D9").Value = "2.604.48"
$ws.Range("E9").Value = "  -2.83%  "
$ws.Range("D10").Value = "'6.72"
$ws.Range("E10").Value = "  +9.21%  "
$ws.Range("E11").Value = "  -1.60%  "
$ws.Range("D12").Value = "'0.347"
$ws.Range("E12").Value = "  -1.22%  "
$ws.Range("E13").Value = "  +1.63%  "
$ws.Range("D14").Value = "3.047.73"
$ws.Range("E14").Value = "  -2.04%  "
$ws.Range("D15").Value = "60.354.13"
$ws.Range("E15").Value = "  -0.87%  "
$ws.Range("D16").Value = "'21.63"
$ws.Range("E16").Value = "  -2.73%  "
$ws.Range("E17").Value = "  -0.50%  "
$ws.Range("D18").Value = "2.596.54"
$ws.Range("E18").Value = "  -2.74%  "
$ws.Range("D19").Value = "'4.76"
$ws.Range("E19").Value = "  -1.32%  "
$ws.Range("D20").Value = "'351.38"
$ws.Range("E20").Value = "  +0.68%  "
$ws.Range("D21").Value = "'10.56"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").Value = "'6.13"
$ws.Range("E22").Value = "  -1.11%  "
$ws.Range("E23").Value = "  +0.27%  "
$ws.Range("D24").Value = "'60.28"
$ws.Range("E24").Value = "  -0.38%  "
$ws.Range("D25").Value = "'0.422"
$ws.Range("E25").Value = "  -0.89%  "
$ws.Range("E26").Value = "  -0.30%  "
$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "  +0.24%  "
$ws.Range("D28").Value = "0.0₃0844"
$ws.Range("E28").Value = "  -3.09%  "
$ws.Range("D29").Value = "'7.37"
$ws.Range("E29").Value = "  -2.80%  "
$ws.Range("E30").Value = "  +0.23%  "
$ws.Range("D31").Value = "'19.42"
$ws.Range("E31").Value = "  -1.22%  "
$ws.Range("D32").Value = "'152.35"
$ws.Range("E32").Value = "  -2.96%  "
$ws.Range("E33").Value = "  -1.57%  "
$ws.Range("D34").Value = "'5.74"
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("D35").Value = "'4.02"
$ws.Range("E35").Value = "  -1.71%  "
$ws.Range("D36").Value = "'1.19"
$ws.Range("E36").Value = "  -3.39%  "
$ws.Range("D37").Value = "'0.863"
$ws.Range("E37").Value = "  +2.74%  "
$ws.Range("E38").Value = "  -3.95%  "
$ws.Range("D39").Value = "'0.845"
$ws.Range("E39").Value = "  -4.86%  "
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").Value = "'3.77"
$ws.Range("E40").Value = "  -0.75%  "
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").Value = "'36.07"
$ws.Range("E41").Value = "  +1.04%  "
$ws.Range("D42").Value = "'298.75"
$ws.Range("E42").Value = "  -3.34%  "
$ws.Range("E43").Value = "  -1.29%  "
$ws.Range("D44").Value = "'0.619"
$ws.Range("E44").Value = "  -4.38%  "
$ws.Range("E45").Value = "  +0.69%  "
$ws.Range("D46").Value = "'0.0554"
$ws.Range("E46").Value = "  -4.28%  "
$ws.Range("D47").Value = "'19.76"
$ws.Range("E47").Value = "  -2.82%  "
$ws.Range("D48").Value = "'4.89"
$ws.Range("E48").Value = "  -3.00%  "
$ws.Range("E49").Value = "  -1.72%  "
$ws.Range("E50").Value = "  +0.16%  "
$ws.Range("D51").Value = "1.996.38"
$ws.Range("E51").Value = "  -2.14%  "
